# Helper: locate $findText anywhere in the document body and replace the
# run(s) that hold it with a hand-built fragment of run/proofErr markup
# (wrapped in a throw-away <w:p> so InsertXML has a valid paragraph-level
# context; since the target Range is narrower than a full paragraph, Word
# merges just the runs into the existing paragraph and keeps its <w:pPr>).
function Replace-WithXml($d, $findText, $innerXml) {
    # Search on a private duplicate of the whole-document range so the
    # Find object's internal state never gets reused for InsertXML below
    # (reusing the same Range object that Find just matched causes the new
    # XML to be inserted rather than substituted for the match).
    $dup = $d.Content.Duplicate
    $found = $dup.Find.Execute($findText, $true, $false, $false, $false, `
                                $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $findText"
        return
    }
    $r = $d.Range($dup.Start, $dup.End)
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:p>'
    $r.InsertXML($xml)
}

$d = $word.ActiveDocument

Replace-WithXml $d 'Training error =  0.02904' '<w:r><w:t xml:space="preserve">Training error </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.02774</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Test error     =  0.0329' '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">Test error     </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>=  0.0328</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Training error =  0.0392' '<w:r><w:t xml:space="preserve">Training error </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.0392</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Test error     =  0.0465' '<w:r><w:t xml:space="preserve">Test error     </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.0465</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Training error =  0.18842' '<w:r><w:t xml:space="preserve">Training error </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.18842</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Test error     =  0.1818' '<w:r><w:t xml:space="preserve">Test error     </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.1818</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Training error =  0.07084' '<w:r><w:t xml:space="preserve">Training error </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.07084</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Test error     =  0.0681' '<w:r><w:t xml:space="preserve">Test error     </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.0681</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Training error =  0.05922' '<w:r><w:t xml:space="preserve">Training error </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.05922</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Test error     =  0.0588' '<w:r><w:t xml:space="preserve">Test error     </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.0588</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Training error =  0.0607' '<w:r><w:t xml:space="preserve">Training error </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.0607</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Test error     =  0.0599' '<w:r><w:t xml:space="preserve">Test error     </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.0599</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Training error =  0.0495' '<w:r><w:t xml:space="preserve">Training error </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.0495</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Test error     =  0.0509' '<w:r><w:t xml:space="preserve">Test error     </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.0509</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Training error =  0.0515' '<w:r><w:t xml:space="preserve">Training error </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.0515</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Test error     =  0.0518' '<w:r><w:t xml:space="preserve">Test error     </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.0518</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Training error =  0.04508' '<w:r><w:t xml:space="preserve">Training error </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">=  </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>0.04508</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Test error     =  0.0466' '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Test error     </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>=  0.0466</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Training error =  0.01434' '<w:r><w:t xml:space="preserve">Training error </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.01434</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Test error     =  0.0324' '<w:r><w:t xml:space="preserve">Test error     </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.0324</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Training error =  0.00264' '<w:r><w:t xml:space="preserve">Training error </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.00264</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Replace-WithXml $d 'Test error     =  0.0293' '<w:r><w:t xml:space="preserve">Test error     </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>=  0.0293</w:t></w:r><w:proofErr w:type="gramEnd"/>'
